$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the accuracy/false-positive/false-negative result cells (C4:F8)
# with the newly measured standard-deviation ("+- ...") annotations.
$ws.Range("C4").Value = "0.903 +-0.023`nfalse pos. 2 false n. 10"
$ws.Range("D4").Value = "0.659 +-0.005`nfalse pos. 0 false n. 43"
$ws.Range("E4").Value = "0.697 +-0.015`nfalse pos. 1 false n. 37"
$ws.Range("F4").Value = "0.653 +-0.002`nfalse pos. 44, false n. 0"

$ws.Range("C5").Value = "0.823 +-0.017`nfalse pos. 0 false n. 22"
$ws.Range("D5").Value = "0.659 +-0.003`nfalse pos. 0 false n. 43"
$ws.Range("E5").Value = "0.697 +- 0.014`nfalse pos. 1 false n. 37"
$ws.Range("F5").Value = "0.653 +-0.002`nfalse pos. 44 false n. 0"

$ws.Range("C6").Value = "0.798 +-0.008`nfalse pos. 0 false n. 26"
$ws.Range("D6").Value = "0.659 +- 0.003`nfalse pos. 0 false n. 43"
$ws.Range("E6").Value = "0.690 +- 0.014`nfalse pos. 1 false n. 39"
$ws.Range("F6").Value = "0.653 +-0.002`nfalse pos. 44 false n. 0"

$ws.Range("C7").Value = "0.8358 +-0.016`nfalse pos. 0 false n. 21"
$ws.Range("D7").Value = "0.659 +- 0.003`nfalse pos. 0 false n. 43"
$ws.Range("E7").Value = "0.697 +- 0.014`nfalse pos. 1 false n. 37"
$ws.Range("F7").Value = "0.653 +-0.002`nfalse pos. 44 false n. 0"

$ws.Range("C8").Value = "0.838 +- 0.015`nfalse pos. 0 false n. 21"
$ws.Range("D8").Value = "0.659 +- 0.003`nfalse pos. 0 false n. 43"
$ws.Range("E8").Value = "0.697 +- 0.014`nfalse pos. 1 false n. 37"
$ws.Range("F8").Value = "0.653 +-0.002`nfalse pos. 44 false n. 0"

# Update the active selection / window position to match the saved view state.
$ws.Range("K8").Select()
$excel.ActiveWindow.Left = 7200
